# fall 22 week 15 complete
# Append the week 15 matchup results (22 rows) to the bottom of Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(4, 1, 4, 2),
    @(2, 1, 2, 2),
    @(7, 0, 5, 2),
    @(7, 2, 6, 0),
    @(3, 0, 5, 3),
    @(5, 2, 5, 1),
    @(6, 2, 5, 0),
    @(2, 2, 5, 0),
    @(6, 2, 4, 1),
    @(4, 0, 4, 2),
    @(4, 3, 4, 0),
    @(4, 1, 4, 2),
    @(5, 0, 5, 2),
    @(3, 1, 2, 2),
    @(4, 0, 3, 3),
    @(5, 0, 4, 2),
    @(4, 0, 7, 3),
    @(6, 0, 5, 2),
    @(6, 2, 6, 1),
    @(3, 2, 6, 0),
    @(5, 2, 4, 1),
    @(5, 0, 2, 2)
)

$startRow = 1503
$endRow = $startRow + $rows.Count - 1

$data = New-Object 'object[,]' $rows.Count,4
for ($i = 0; $i -lt $rows.Count; $i++) {
    for ($j = 0; $j -lt 4; $j++) {
        $data[$i, $j] = $rows[$i][$j]
    }
}

$rng = $ws.Range("A$startRow`:D$endRow")
$rng.Value = $data

# Mirror the view state recorded in the saved file: scrolled near the
# newly-added rows with the next empty row selected.
$newSelRow = $endRow + 1
$ws.Range("A$newSelRow").Select()
$excel.ActiveWindow.ScrollRow = $startRow + 9
$excel.ActiveWindow.ScrollColumn = 1
